# Daily attendance processing - 2025-11-02 09:20:21
# Normalizes the "Recorded By" (column G) lists on the "Session Analysis
# Results" sheet: the first name/e-mail in each comma-separated list is
# moved to the end of the list (left rotation by one), except for the
# literal value "System, admin@admin.com" which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($text -eq $null) { continue }
    if ($text -eq "") { continue }
    if ($text -eq "System, admin@admin.com") { continue }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
